$d = $word.ActiveDocument

$d.TrackRevisions = $true

$searchRange = $d.Content
$searchRange.Find.ClearFormatting()
while ($searchRange.Find.Execute("4", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $searchRange.Text = "6"
    $searchRange.Collapse(0)
    $searchRange.End = $d.Content.End
}

$d.TrackRevisions = $false

while ($d.Revisions.Count -gt 0) {
    $d.Revisions.Item(1).Accept()
}
